$d = $word.ActiveDocument

# Remove the "Boekingsadres" and "Boekingsplaats" paragraphs entirely,
# including their paragraph marks, leaving "Boekingsdatum" as-is.
$targets = @("Boekingsadres", "Boekingsplaats")

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    $trimmed = $text.Trim()
    if ($targets -contains $trimmed) {
        $para.Range.Delete()
    }
}
